$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values
$ws.Range("B2").Value = 5.6
$ws.Range("B3").Value = 5.6
$ws.Range("C4").Value = 1.5

# Remove the "theta_threshold_range" row (row 5), shifting the row below
# ("pie_threshold_range") up into row 5.
$ws.Rows.Item(5).Delete()

# Update the (now shifted-up) pie_threshold_range row values
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Adjust column C width to match the new best-fit width
$ws.Columns.Item(3).ColumnWidth = 5.5

# Set page setup (paper size / orientation) to match the saved print settings
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
